$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 47, shifting existing rows 47:170 down to 48:171.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new weekly record.
$ws.Cells.Item(47, 1).Value = 8
$ws.Cells.Item(47, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 44497
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = 100112012
$ws.Cells.Item(47, 7).Value = "Espinaca"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 3000
$ws.Cells.Item(47, 11).Value = 400
$ws.Cells.Item(47, 12).Value = 500
$ws.Cells.Item(47, 13).Value = 450
$ws.Cells.Item(47, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(47, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(47, 16).Value = 900
$ws.Cells.Item(47, 17).Value = 0.5
$ws.Cells.Item(47, 18).Value = "Hortaliza"
